$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update the series title in C1 (dropped "SAR (China)")
$ws.Range("C1").Value = "(DC)Hong Kong Retail Bonds: Price: Mid: HK Link A: 07-05-2009: 3.60%"

# 3. Relabel the A11 caption
$ws.Range("A11").Value = "Function Information"

# 4. Correct the "Last Update Time" for the second series (C14)
$ws.Range("C14").Value = 41781

# 5/6. Tiny floating point corrections to the Skewness/Kurtosis stats
$ws.Range("B20").Value = -0.2234505289774544
$ws.Range("B21").Value = -1.733729408587631

# 7. Widen the custom numeric format used by the historical data rows
$ws.Range("B27:C46").NumberFormat = "###0.000"

# 8. Refresh the embedded CEIC add-in metadata comment on A1 to the
#    regenerated payload (series name updated, graph title populated, etc.)
$comment = $ws.Range("A1").Comment
$comment.Text("WCkAAB+LCAAAAAAAAAPtWltvG8cV/isLAgYSINTukpJFEeMNeJFswqQkiHRt5yVY7o7EqZa77M6sJL65KAK7bRKgcJPUzaVNkbRNUdsBnKZu7ebHBCYtPfkv9MzMXknKJh0HaQMbhLxzbnPmzJlzvr2g14/6jnKAfUo891xOX9JyCnYtzybu3rlcwHbz+tnc6wZaP7Kws236Zh8zEFZAy6XlI0rO5XqMDcqqenh4uHRYXPL8PbWgabp6pdVsWz3cN/PEpcx0LZyLtexna+UMVLP7LcxM22Sm1DyXa7QbSzVMrDrQWqZr7mF/qRpQ4mJK111GGMGUa/rYZLhWb/1ELswoLJ1d0pE6RU8kqwFxbCmXkZT0UA6mxR3Sx0ZB00t5bS2vFTu6Xi5oZW11qVDU3ogUY0HUNClrY/+AWILQZmZ/INS1Na2o67DoVaTOFAJbSQAMtOXYO/iAUGzXsOPQhSKihhtYsRiserFgakhN6YaGnt+F87456HUIc/BibpxRBthXWp7Leq8pFy7WQ68Scwba8HxsQSSfy7lNfLjlhwHuDJrA7fSIz4Z1c7iwrUsU+1sDHq7FVA1Uh/VVHOyzSwPYdWxDUgDDYH6AkXoKM1GqE2rBNXEDbBu7pkPTShkmuuz5+3RgWngTTrTKbRy6jmfakHqMUEasZNIpBtr2PdgM7nLVc+wNsBoKz2DElhsuhJhPW/W8/cS7WUwkdlXsL+xp32SR+BQdtXve4ZbrDNtBl1o+6WK7Xo2kZ/IQP5qhdi2gzOuDFwkJSVqKYttqq6UO4R8cyEkmqmOL9E1n24FQUqMItjIEVAmYt0tYzXOCvksjzyao6DKsq4OP4nXGY7QFW+zy0Htuw43kZbBnsrIKO95hPOc0Q4QiRa5QK9r0acakcB1o0SZOc8S+8FVuEAf6RXpHUtRsbrR7GLOZiSE5iJfGDd6BjOpwM+h34ZB14aQdiFkpUhM+gmyFjAe/DA26Sl78OhqUaf4DP2I2Wnft0+UiJoLpUnMZOvAmSAjW5FQd090H6mXCepuVaC0zOEhG4FT5aR6C8ztwzKEgx1FK01DDtZzAxrIsNNxdkaLcN7mpp7LRFKkJJ91ApjvsDAdQpykpM7g4l4POXabMB2yQMywvcJk/5PUDqaHos3Ro0HXFBKYzt86uj38WACQZbgSuVfPs+WezZXQuuYTN76EX+LIozq8iosfrY0DrmFcaUfrn1rcWWRP1FxLvu7jvucSaP9oQZO69/RwLodGpmlsDy/M1t7wD3V02P37W51bzAU9Cu1tomgqlnkVEsobHw07pq6ccmTreNQMHsByDRrsX195JMqrQ/UmZNAld8p2oAhocKVOAypbdX7IAQnA4uGR5fU5QAaFebiM1Lc9xkIXX3b2m6e4FgDTiujJJj+sv75Id33QpX04MLCZK8WwhFNUpCXgMWby2ApEIsnh5wEXqhBzq4P7A802nBYEhG2HahagJAEnLZL1wBL3NwVYUZDVRjbWynkWOP0tMNCm5DH7gwzI5QRRCfC0SlCcyCQ3xVbbgWDo10yFdX1bVqJXP4sGGJRAxqr98cQvCxWgP4MYMuu9FPORoPRmEdJGyesSQCcwLqdHeWS4VltdKK4Bs+BiJFdfxwKOEKTuwwLJS6xGA6jVQhEgqVdjDslI52CsrG8SFmzoCxPaQMtwvK0VNYZ5SWlM4ZoYGKcIXtqgXajltE21E/UEsoQ5ryVKyAoBT9sC8My0YcxIFo24SZ5iSkxFqehaIjW7/bvz3b0Y3vnj87vUnD38/vvVvuH50//boszsnN78+/vRtII7+dW/0MeeeXP/N8SfvPPrPzcf3HsCwqB1fv1daG332hQyTtIg6ZtfBwulOtVTSisuQwzEJ8Y1TBfi2A4sJ2tWrAnPHYxTeRIpBbb1RO9+siloVEyN12a5Ufn869IJk2JYrFROJdFGjLJMiRieqfeE4w021P4PfLx7grHSaf5qijMXjB395/OD2qdphwBIcp6+treQ1+BWeDvPgznuZy+lrs2BeM9NfYuFC2uiEDNqRzSWOU8M2ijrc3K+sanrcJ+z4IMwSmmSFljrmnjqhJ0k1Cb3iFEiPI6Y4HB1MWcyWxyU1kEEc/+PXx3ffy0iF0Q0pWSvgnIBIfDI1GgjTmzsdpb11aae2rnTW2zxPEl5KThp/inA4e3zo0kkjDuRrSm5F8RLsnUud0EyGzaJK64vYmfTnvO8FA1lsUgoJdYZkUnVmacyoPIIna+ZECUo4M6TDHf3g81kK4TrqCSpOP1SJs1Aw0DRb0lNC4WQf3ZA18EzGQjhZ/JgBshrOTnoYJzkUuLBzTVDQ5bYI6L72ZqpThUR+S7btEZdRY3l1RRf3Y+EYgbLO7Yn/UaMP7VOYFmED+gQFXTDp+hELD7KxidQsATwdmNC5veQ+NibImp2E9+TjT8YffjV+/8vj638b/eqvo7fff/zgD8d3PpWnbPzel+N37oRVfbLwC1/43bEElIp44mIp/PQpHAco3167qbgeUwC+KIGoQN9eu5Uyxh0VQCexDPAwdiTrwpRoWpnrKSlXYh8yerGKBBM13rKKsUTYtLwBsZJJ3shzU4BxFcF4pdHJBxQrHiCzV2ElWeFEeV69UEW20O1VraAXQq70hi+ha9JU6M87XheARsQQDzMmRDJaT1dIZMV855tb1UozEZFObPk29nkaygsUwVPeQho0GkWplqIAF0CkFTj8+dOU2DQrtpwqZmr4IGe3YvNH/LMfe2QkUC3wfQmS3PA9QTsYALCOn9WcLiCehKbA9KYEvml4nYwb9SwfxikudL4smxMEXxSmkCWLVIPyh0YSG2/y2CRD4GWenkI8wncBElodENgTlZeedd/3/Jn1J+FEYi2A5VBS1CTksYzYVAnh7WSzIkJU814UlC+WMlD+lXrt1QsenJqL/M8OZtBMlCrc3dOysu0DJiwrLWKXlQsXlSZx95VKWdFWBeLRtDWA3UtntTOzgPyLsfs/AeOfPLwx+vlH49t/kmX7ycNfnvz51vj+3ZMPvx799s7oGqD6+4DYHz345/iPD/nF/dsnH3w1usavQ8mb34zf/bwC42SNHOMnwXsJ778rvIeYArYHML78LHjPX8qB6NnvHd6X5oH3pZfw/iW8/8HgvXhNOg3rBfk0OH//7uitX3xHKK+fDuX1N1OdahrK69rZ0kso/xLK/19Bef1HAOXlc+0fPZKfKD/fC5JXk8f5cWX06tjBbLFPXtREu+UdPLcu7P2iqg265dhhMBd7IxGHJTGQ/u6HJ8qL++xHpl3F9wE+8S8EFv5OJ3qztWO6ewv6JRclFPkrIpg9/JBmg/iUXeGtKbySlKsx5arEoVeMggSbV+T4qrG8IgkgoKatqxk3o0PM5HdYntMkfbLgeyMtOulZIxDLwUACtcZiOcObzCY+AhiZsgDlsftTaCDyResi1mTqQlWN9fnHDZTs9diijq12TWzjrpa3uriQX7a1Un4N42Je1+GvaRX4jQX/NCI0DjWE4MMFJ1GjDUu+jzT+C+biD6VYKQAA")
